# "Updated split payment test data"
# The custcredit value (row 2, column X) changes from "100" to "0".
# The original value was stored as text (quote-prefixed) rather than a
# number, so we use a leading apostrophe to force the same text storage
# and keep the cell's existing "quotePrefix" style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("X2").Value = "'0"

# The user also scrolled the sheet right and re-selected a cell (X3) before
# saving - reflect that view/selection state too.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 18
$ws.Range("X3").Select()
